$d = $word.ActiveDocument

# Locate the list-item paragraph that currently starts with "Search " --
# this is the paragraph the commit rewrites.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.StartsWith("Search ")) {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $apos = [char]0x2019
    $newLead = "Let" + $apos + "s play Photoshop"

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
           'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
           'w14:paraId="11C39852" w14:textId="002E691A" w:rsidR="00ED4143" w:rsidRPr="009C6E70" ' +
           'w:rsidRDefault="000D50D9" w:rsidP="000D50D9">' +
             '<w:pPr>' +
               '<w:pStyle w:val="ListParagraph"/>' +
               '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
               '<w:jc w:val="both"/>' +
               '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' +
             '</w:pPr>' +
             '<w:r w:rsidRPr="009C6E70">' +
               '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' +
               '<w:t>' + $newLead + '</w:t>' +
             '</w:r>' +
             '<w:r w:rsidR="002370BB" w:rsidRPr="009C6E70">' +
               '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' +
               '<w:t>.</w:t>' +
             '</w:r>' +
             '<w:r w:rsidRPr="009C6E70">' +
               '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/></w:rPr>' +
               '<w:t xml:space="preserve"> Using the Lena image, apply two different mathematical algorithms/equations. Explain what happened in terms of numbers and visualization.</w:t>' +
             '</w:r>' +
           '</w:p>'

    $target.Range.InsertXML($xml) | Out-Null
}
